# Active_Outages.xlsx refresh — 6/19/2025, 9:13:13 AM
# Updates the "Elapsed Duration(Hrs)" column on several existing outage rows
# (time has moved on since the report was last generated) and appends a
# brand-new outage row to the R1 sheet.

$wb = $excel.ActiveWorkbook

# ---- R1 : elapsed-duration refresh ------------------------------------
$ws1 = $wb.Worksheets.Item("R1")
$ws1.Cells.Item(2, 7).Value = "3946:27:15"
$ws1.Cells.Item(3, 7).Value = "85:59:53"
$ws1.Cells.Item(4, 7).Value = "108:59:53"

# ---- R1 : new outage row (row 6) ---------------------------------------
$ws1.Cells.Item(6, 1).Value  = ""
$ws1.Cells.Item(6, 2).Value  = "R4"
$ws1.Cells.Item(6, 3).Value  = ""
$ws1.Cells.Item(6, 4).Value  = "JED0125"
$ws1.Cells.Item(6, 5).Value  = ""
$ws1.Cells.Item(6, 6).Value  = "2025-06-19 09:13:02"
$ws1.Cells.Item(6, 7).Value  = "0:00:00"
$ws1.Cells.Item(6, 8).Value  = ""
$ws1.Cells.Item(6, 9).Value  = "Generator-SG"
$ws1.Cells.Item(6, 10).Value = "Good+In progress"
$ws1.Cells.Item(6, 11).Value = ""
$ws1.Cells.Item(6, 12).Value = "Latis"

# ---- R2 : elapsed-duration refresh -------------------------------------
$ws2 = $wb.Worksheets.Item("R2")
$ws2.Cells.Item(2, 7).Value = "12127:49:56"
$ws2.Cells.Item(3, 7).Value = "3257:33:25"
$ws2.Cells.Item(4, 7).Value = "495:44:59"

# ---- R4 : elapsed-duration refresh -------------------------------------
$ws4 = $wb.Worksheets.Item("R4")
$ws4.Cells.Item(2, 7).Value = "2973:39:45"
$ws4.Cells.Item(3, 7).Value = "200:52:00"
$ws4.Cells.Item(4, 7).Value = "89:04:25"
$ws4.Cells.Item(5, 7).Value = "86:41:58"

# ---- R5 : elapsed-duration refresh -------------------------------------
$ws5 = $wb.Worksheets.Item("R5")
$ws5.Cells.Item(2, 7).Value = "447:38:44"

# ---- R6 : elapsed-duration refresh -------------------------------------
$ws6 = $wb.Worksheets.Item("R6")
$ws6.Cells.Item(2, 7).Value = "88:11:02"
